# Adds 5 new price-report rows (443-447) for "Comercializadora del Agro de
# Limarí" - Naranja, week of 2022-09-28 (serial date 44832).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ K = "Lane Late";  L = "Especial"; N = 100000; O = 110000; P = 105000; S = 262 },
    @{ K = "Lane Late";  L = "Primera";  N = 80000;  O = 90000;  P = 85000;  S = 212 },
    @{ K = "Lane Late";  L = "Segunda";  N = 50000;  O = 60000;  P = 55000;  S = 138 },
    @{ K = "Navel Late"; L = "Primera";  N = 80000;  O = 90000;  P = 85000;  S = 212 },
    @{ K = "Navel Late"; L = "Segunda";  N = 50000;  O = 60000;  P = 55000;  S = 138 }
)

$startRow = 443
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = 2
    $ws.Cells.Item($r, 2).Value = "Comercializadora del Agro de Limarí"
    $ws.Cells.Item($r, 3).Value = "Coquimbo"
    $ws.Cells.Item($r, 4).Value = 44832
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 5).Value = 4
    $ws.Cells.Item($r, 6).Value = "Fruta"
    $ws.Cells.Item($r, 7).Value = 100102
    $ws.Cells.Item($r, 8).Value = "Cítricos"
    $ws.Cells.Item($r, 9).Value = 100102005
    $ws.Cells.Item($r, 10).Value = "Naranja"
    $ws.Cells.Item($r, 11).Value = $data.K
    $ws.Cells.Item($r, 12).Value = $data.L
    $ws.Cells.Item($r, 13).Value = 20
    $ws.Cells.Item($r, 14).Value = $data.N
    $ws.Cells.Item($r, 15).Value = $data.O
    $ws.Cells.Item($r, 16).Value = $data.P
    $ws.Cells.Item($r, 17).Value = "$/bins (400 kilos)"
    $ws.Cells.Item($r, 18).Value = "Provincia de Limarí"
    $ws.Cells.Item($r, 19).Value = $data.S
    $ws.Cells.Item($r, 20).Value = 400
}
